# Automatische test-sync: 2025-08-13 20:57:50
#
# Appends three new "Demo inplannen" log rows (rows 6-8) to the Logs
# sheet - mirroring the existing rows 2-5 - and refreshes the Dashboard
# summary count to match the new total.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newTimestamps = @(
    "2025-08-13 20:57:38",
    "2025-08-13 20:57:40",
    "2025-08-13 20:57:41"
)

$row = 6
foreach ($ts in $newTimestamps) {
    $logs.Cells.Item($row, 1).Value = "Demo inplannen"
    $logs.Cells.Item($row, 2).Value = "klantenservice@testbedrijf123.nl"
    $logs.Cells.Item($row, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
    $logs.Cells.Item($row, 4).Value = "Intern verzoek / Actie voor medewerker"
    $logs.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
    $logs.Cells.Item($row, 6).Value = $ts
    $logs.Cells.Item($row, 7).Value = "Nee"
    $logs.Cells.Item($row, 8).Value = "Ja"
    $logs.Cells.Item($row, 9).Value = "Nee"
    $logs.Cells.Item($row, 10).Value = "Nee"
    $row++
}

# Keep the per-category conditional formatting ranges in sync with the
# new bottom row of data (D/G/H/I/J previously stopped at row 5). Use
# ModifyAppliesToRange on each existing rule so priority/dxfId/formula
# stay untouched - only the applied range grows.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$($col)2:$($col)5")
    $newRange = $logs.Range("$($col)2:$($col)8")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Dashboard summary: 4 existing "Intern verzoek" rows + 3 new ones = 7.
$dashboard.Range("B2").Value = 7
